$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 9 (Ano 2025) with refreshed figures
$ws.Range("B9").Value = 2932227.58
$ws.Range("C9").Value = 466376.99
$ws.Range("D9").Value = 3398604.57
$ws.Range("E9").Value = 13.72260233263913
$ws.Range("F9").Value = 86.27739766736087
$ws.Range("G9").Value = -54.92698937937983
$ws.Range("H9").Value = -47.0480229463199
$ws.Range("I9").Value = 29324
$ws.Range("J9").Value = 1253
$ws.Range("K9").Value = 30577
$ws.Range("L9").Value = 21094
$ws.Range("M9").Value = 161.1171219304068
$ws.Range("N9").Value = 9.997940731831466
